$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add name + attendance status in the previously-empty A2/B2 cells
$ws.Range("A2").Value = "богдана"
$ws.Range("B2").Value = "Відсутній"

# Row 3: summary count for "Присутність" (attendance) column goes from 0 to 1
$ws.Range("B3").Value = 1

# Row 5: "Будуть" (will attend) count and "Відмітилось" (checked-in) count go from 0 to 1
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 1
